$p = $ppt.ActivePresentation

# --- Slide 12 ("Next Steps"): Content Placeholder 2, second paragraph ---
# Before: "Request SPRING " + "WG adoption"  (2 runs)
# After:  "In " + "SPRING " + "WG adoption poll"  (3 runs)
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange
$para2 = $tr12.Paragraphs(2)

# Drop the leading "Request " (8 characters) from the run that starts the paragraph,
# leaving behind "SPRING " in what used to be the first run.
$prefix = $para2.Characters(1, 8)
[void]$prefix.Delete()

# Add the new leading text "In " as its own run at the start of the paragraph.
[void]$para2.InsertBefore("In ")

# Re-anchor on the (shifted) paragraph/run and extend "WG adoption" -> "WG adoption poll".
$para2 = $tr12.Paragraphs(2)
$lastRun = $para2.Runs(3, 1)
$lastRun.Text = "WG adoption poll"

# --- Slide 3 ("Requirements and Scope"): Content Placeholder 2, second paragraph ---
# "Delay and Loss Performance Measurement (PM) " -> "Delay and Synthetic Loss Performance Measurement (PM) "
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$reqPara = $tr3.Paragraphs(2)
$reqRun = $reqPara.Runs(1, 1)
$reqRun.Text = "Delay and Synthetic Loss Performance Measurement (PM) "
